$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.893.65"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.635.57"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").Value = "215.43"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").Value = "28.78"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "0.260"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "1.869.47"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "1.635.36"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "0.586"
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").Value = "9.45"
$ws.Range("E15").Value = "  +6.12%  "
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "29.907.05"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "64.76"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "240.16"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "9.92"
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  +4.05%  "
$ws.Range("D25").Value = "157.46"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "15.53"
$ws.Range("D27").Value = "0.109"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "6.64"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "3.20"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "1.422.72"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("D37").Value = "2.76"
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").Value = "76.01"
$ws.Range("E40").Value = "  +9.64%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "0.0500"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").Value = "1.776.92"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "48.74"
$ws.Range("E49").Value = "  -9.54%  "
$ws.Range("D50").Value = "93.09"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("E51").Value = "  +7.53%  "
